# Allgemeine Aktualisierung: neue Tageswerte (Montag 2020-11-09 bis Freitag
# 2020-12-04) an die bestehende "Rohdaten"-Tabelle anhängen.
# Spalten: A=Datum, B=Uhrzeit, C=Summe, D=Summe Aachen, E=Summe Todesfälle,
#          F=Summe genesen, G=Akute Fälle

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Datum (serial), Uhrzeit (text), Summe, Summe Aachen,
# Summe Todesfälle, Summe genesen, Akute Fälle
$newRows = @(
    @(44144, "10:30:00", 7273, 3193, 132, 5831, 1310),
    @(44145, "10:00:00", 7426, 3248, 140, 6076, 1210),
    @(44146, "10:45:00", 7658, 3335, 142, 6237, 1279),
    @(44147, "09:15:00", 7871, 3417, 143, 6394, 1334),
    @(44148, "10:45:00", 8211, 3552, 143, 6638, 1430),
    @(44151, "09:45:00", 8681, 3738, 143, 7124, 1414),
    @(44152, "09:20:00", 8728, 3768, 148, 7322, 1258),
    @(44153, "10:30:00", 8920, 3840, 148, 7498, 1274),
    @(44154, "09:20:00", 9179, 3950, 156, 7631, 1392),
    @(44155, "09:30:00", 9311, 3995, 162, 7849, 1300),
    @(44158, "10:30:00", 9565, 4100, 164, 8230, 1171),
    @(44159, "09:30:00", 9676, 4159, 168, 8447, 1061),
    @(44160, "10:20:00", 9827, 4226, 172, 8576, 1079),
    @(44161, "10:30:00", 9867, 4241, 175, 8698,  994),
    @(44162, "10:20:00", 10038, 4325, 178, 8786, 1074),
    @(44165, "10:00:00", 10289, 4425, 182, 8868, 1239),
    @(44166, "11:15:00", 10316, 4441, 184, 8968, 1164),
    @(44167, "10:15:00", 10460, 4487, 188, 9194, 1078),
    @(44168, "10:00:00", 10605, 4550, 191, 9314, 1100),
    @(44169, "10:00:00", 10764, 4627, 196, 9548, 1020)
)

# First new row is 177 - one past the current last row (176)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$row = $lastRow + 1

foreach ($entry in $newRows) {
    # Copy the date cell above so the new row picks up the same date style
    # (custom "YYYY-MM-DD HH:MM:SS" number format with border/centering)
    # instead of Excel inventing a brand-new style index for it.
    $ws.Cells.Item($row - 1, 1).Copy($ws.Cells.Item($row, 1))

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]

    $row = $row + 1
}
